$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F2").Value = 760
$ws.Range("F3").Value = 77
$ws.Range("F4").Value = 2047
$ws.Range("F5").Value = 6109
$ws.Range("F6").Value = 3506
$ws.Range("G6").Value = 80
$ws.Range("F8").Value = 57
$ws.Range("F9").Value = 1457
$ws.Range("F10").Value = 4854
$ws.Range("G10").Value = 68
$ws.Range("F11").Value = 1118
$ws.Range("F12").Value = 1809
$ws.Range("F13").Value = 27
$ws.Range("F14").Value = 70
$ws.Range("F15").Value = 76
$ws.Range("F16").Value = 224
$ws.Range("F17").Value = 167
$ws.Range("F18").Value = 1071
$ws.Range("F19").Value = 346
$ws.Range("F22").Value = 99
$ws.Range("F23").Value = 11
$ws.Range("F26").Value = 41
$ws.Range("F27").Value = 1166
$ws.Range("F28").Value = 447
$ws.Range("F29").Value = 129
$ws.Range("F30").Value = 250
$ws.Range("F31").Value = 507
$ws.Range("F34").Value = 1869
$ws.Range("F35").Value = 2328
$ws.Range("F36").Value = 1100
$ws.Range("F38").Value = 42
$ws.Range("F39").Value = 305
$ws.Range("F40").Value = 146
$ws.Range("F41").Value = 683
$ws.Range("F42").Value = 536
$ws.Range("F43").Value = 72
$ws.Range("F44").Value = 699
$ws.Range("F45").Value = 71
$ws.Range("F46").Value = 489
$ws.Range("F47").Value = 529
$ws.Range("F49").Value = 168

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F9").Value = 164
$ws.Range("F16").Value = 126

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F3").Value = 760
$ws.Range("F4").Value = 77
$ws.Range("F5").Value = 2047
$ws.Range("F6").Value = 3506
$ws.Range("G6").Value = 80
$ws.Range("F7").Value = 57
$ws.Range("F8").Value = 1457
$ws.Range("F9").Value = 4854
$ws.Range("G9").Value = 68
$ws.Range("F10").Value = 1809
$ws.Range("F11").Value = 27
$ws.Range("F13").Value = 70
$ws.Range("F16").Value = 76
$ws.Range("F17").Value = 224
$ws.Range("F18").Value = 167
$ws.Range("F19").Value = 164
$ws.Range("F20").Value = 1071
$ws.Range("F21").Value = 346
$ws.Range("F22").Value = 99
$ws.Range("F25").Value = 41
$ws.Range("F26").Value = 1166
$ws.Range("F27").Value = 447
$ws.Range("F28").Value = 129
$ws.Range("F29").Value = 250
$ws.Range("F32").Value = 1869
$ws.Range("F33").Value = 2328
$ws.Range("F34").Value = 1100
$ws.Range("F36").Value = 42
$ws.Range("F37").Value = 305
$ws.Range("F38").Value = 146
$ws.Range("F41").Value = 683
$ws.Range("F42").Value = 536
$ws.Range("F43").Value = 699
$ws.Range("F44").Value = 489
$ws.Range("F45").Value = 529
$ws.Range("F48").Value = 168
